# Weekly update: insert a new price record as the first row of this
# market/product block (row 125), pushing the existing rows down by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(125).Insert()

$ws.Cells.Item(125, 1).Value = 4
$ws.Cells.Item(125, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(125, 3).Value = "Los Lagos"
$ws.Cells.Item(125, 4).Value = 44504
$ws.Cells.Item(125, 5).Value = 10
$ws.Cells.Item(125, 6).Value = "Fruta"
$ws.Cells.Item(125, 7).Value = 100101
$ws.Cells.Item(125, 8).Value = "Berries"
$ws.Cells.Item(125, 9).Value = 100101007
$ws.Cells.Item(125, 10).Value = "Kiwi"
$ws.Cells.Item(125, 11).Value = "Hayward"
$ws.Cells.Item(125, 12).Value = "Primera"
$ws.Cells.Item(125, 13).Value = 200
$ws.Cells.Item(125, 14).Value = 16000
$ws.Cells.Item(125, 15).Value = 17000
$ws.Cells.Item(125, 16).Value = 16500
$ws.Cells.Item(125, 17).Value = "$/caja 15 kilos"
$ws.Cells.Item(125, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(125, 19).Value = 1100
$ws.Cells.Item(125, 20).Value = 15
